$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Range("H9").Value = 333.33334
$ws.Range("I9").Value = 375
$ws.Range("J9").Value = 250
$ws.Range("K9").Value = 375
$ws.Range("L9").Value = 250
$ws.Range("M9").Value = -206
$ws.Range("N9").Value = -588
$ws.Range("H12").Value = 397.57144
$ws.Range("I12").Value = 397.57144
$ws.Range("K12").Value = 397.57144
$ws.Range("M12").Value = -227.57144
$ws.Range("H19").Value = 9995
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 9995
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 9995
$ws.Range("M19").Value = ""
$ws.Range("N19").Value = -10345
$ws.Range("H32").Value = 4235
$ws.Range("J32").Value = 5043.75
$ws.Range("L32").Value = 5043.75
$ws.Range("N32").Value = -5695.75
$ws.Range("I33").Value = 120.07692
$ws.Range("J33").Value = 333333380
$ws.Range("K33").Value = 120.07692
$ws.Range("L33").Value = 333333380
$ws.Range("M33").Value = 108.92308
$ws.Range("N33").Value = -333333838
$ws.Range("H38").Value = 8354.883
$ws.Range("J38").Value = 11082.833
$ws.Range("L38").Value = 33248.499
$ws.Range("N38").Value = -33992.499
$ws.Range("H43").Value = 13424.75
$ws.Range("I43").Value = 16499.5
$ws.Range("K43").Value = 16499.5
$ws.Range("M43").Value = -16430.5
$ws.Range("H49").Value = 6899.8
$ws.Range("J49").Value = 7499.5
$ws.Range("L49").Value = 22498.5
$ws.Range("N49").Value = -22770.5
$ws.Range("H52").Value = 779.625
$ws.Range("J52").Value = 2499
$ws.Range("L52").Value = 7497
$ws.Range("N52").Value = -7817
$ws.Range("H59").Value = 6249.5
$ws.Range("J59").Value = 6249.5
$ws.Range("L59").Value = 18748.5
$ws.Range("N59").Value = -19862.5
$ws.Range("H62").Value = 5333.3335
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 5333.3335
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""
$ws.Range("H86").Value = 1906986.6
$ws.Range("I86").Value = 2354.5557
$ws.Range("K86").Value = 2354.5557
$ws.Range("M86").Value = -1231.5557
$ws.Range("H88").Value = 1436271.2
$ws.Range("I88").Value = 8000
$ws.Range("J88").Value = 1674316.5
$ws.Range("K88").Value = 8000
$ws.Range("L88").Value = 1674316.5
$ws.Range("M88").Value = -7594
$ws.Range("N88").Value = -1675128.5
$ws.Range("H89").Value = 1906986.6
$ws.Range("I89").Value = 2354.5557
$ws.Range("K89").Value = 11772.7785
$ws.Range("M89").Value = -6156.7785
$ws.Range("H91").Value = 1436271.2
$ws.Range("I91").Value = 8000
$ws.Range("J91").Value = 1674316.5
$ws.Range("K91").Value = 8000
$ws.Range("L91").Value = 1674316.5
$ws.Range("M91").Value = -6596
$ws.Range("N91").Value = -1677124.5
$ws.Range("H98").Value = 2739
$ws.Range("I98").Value = 3112.2856
$ws.Range("J98").Value = 997
$ws.Range("K98").Value = 3112.2856
$ws.Range("L98").Value = 997
$ws.Range("M98").Value = -1614.2856
$ws.Range("N98").Value = -3993
$ws.Range("H112").Value = 1070.96
$ws.Range("J112").Value = 1149.381
$ws.Range("L112").Value = 3448.143
$ws.Range("N112").Value = -5664.143
$ws.Range("H113").Value = 11835.429
$ws.Range("I113").Value = 12500
$ws.Range("J113").Value = 11337
$ws.Range("K113").Value = 12500
$ws.Range("L113").Value = 11337
$ws.Range("M113").Value = -9246
$ws.Range("N113").Value = -17845
$ws.Range("H116").Value = 5535.2856
$ws.Range("I116").Value = 6282.6665
$ws.Range("J116").Value = 4974.75
$ws.Range("K116").Value = 6282.6665
$ws.Range("L116").Value = 4974.75
$ws.Range("M116").Value = -2840.6665
$ws.Range("N116").Value = -11858.75
$ws.Range("H122").Value = 2739
$ws.Range("I122").Value = 3112.2856
$ws.Range("J122").Value = 997
$ws.Range("K122").Value = 9336.856800000001
$ws.Range("L122").Value = 2991
$ws.Range("M122").Value = -6886.856800000001
$ws.Range("N122").Value = -7891
$ws.Range("H129").Value = 1869.091
$ws.Range("I129").Value = 1538
$ws.Range("J129").Value = 2145
$ws.Range("K129").Value = 4614
$ws.Range("L129").Value = 6435
$ws.Range("M129").Value = 386
$ws.Range("N129").Value = -16435
$ws.Range("H131").Value = 4999.5
$ws.Range("J131").Value = 4999
$ws.Range("L131").Value = 14997
$ws.Range("N131").Value = -25077
$ws.Range("H132").Value = 2650
$ws.Range("I132").Value = 2650
$ws.Range("K132").Value = 7950
$ws.Range("M132").Value = -5420
$ws.Range("H137").Value = 3827.5
$ws.Range("I137").Value = 4088.5715
$ws.Range("K137").Value = 12265.7145
$ws.Range("M137").Value = -9715.7145
$ws.Range("H138").Value = 1778.2593
$ws.Range("J138").Value = 2470.6453
$ws.Range("L138").Value = 7411.9359
$ws.Range("N138").Value = -17691.9359
$ws.Range("H141").Value = 4429.185
$ws.Range("I141").Value = 3590.8696
$ws.Range("J141").Value = 9249.5
$ws.Range("K141").Value = 10772.6088
$ws.Range("L141").Value = 27748.5
$ws.Range("M141").Value = -5592.6088
$ws.Range("N141").Value = -38108.5
$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Range("H2").Value = 2095.15
$ws.Range("I2").Value = 1167
$ws.Range("K2").Value = 1167
$ws.Range("M2").Value = -1054
$ws.Range("H32").Value = 2829.4565
$ws.Range("I32").Value = 1420.8206
$ws.Range("K32").Value = 1420.8206
$ws.Range("M32").Value = -1133.8206
$ws.Range("H44").Value = 51990
$ws.Range("J44").Value = 51990
$ws.Range("L44").Value = 51990
$ws.Range("N44").Value = -52966
$ws.Range("H45").Value = 2273.5334
$ws.Range("I45").Value = 2009.5454
$ws.Range("J45").Value = 2999.5
$ws.Range("K45").Value = 2009.5454
$ws.Range("L45").Value = 2999.5
$ws.Range("M45").Value = -1632.5454
$ws.Range("N45").Value = -3753.5
$ws.Range("H61").Value = 2698.5
$ws.Range("I61").Value = 2698.5
$ws.Range("K61").Value = 2698.5
$ws.Range("M61").Value = -2486.5
$ws.Range("H63").Value = 2897
$ws.Range("I63").Value = 3586.5557
$ws.Range("J63").Value = 1862.6666
$ws.Range("K63").Value = 3586.5557
$ws.Range("L63").Value = 1862.6666
$ws.Range("M63").Value = -2900.5557
$ws.Range("N63").Value = -3234.6666
$ws.Range("H66").Value = 2897
$ws.Range("I66").Value = 3586.5557
$ws.Range("J66").Value = 1862.6666
$ws.Range("K66").Value = 17932.7785
$ws.Range("L66").Value = 9313.333000000001
$ws.Range("M66").Value = -14500.7785
$ws.Range("N66").Value = -16177.333
$ws.Range("H74").Value = 3296.9375
$ws.Range("I74").Value = 3366.7693
$ws.Range("J74").Value = 2994.3333
$ws.Range("K74").Value = 3366.7693
$ws.Range("L74").Value = 2994.3333
$ws.Range("M74").Value = -2492.7693
$ws.Range("N74").Value = -4742.3333
$ws.Range("H77").Value = 3296.9375
$ws.Range("I77").Value = 3366.7693
$ws.Range("J77").Value = 2994.3333
$ws.Range("K77").Value = 16833.8465
$ws.Range("L77").Value = 14971.6665
$ws.Range("M77").Value = -12465.8465
$ws.Range("N77").Value = -23707.6665
$ws.Range("H116").Value = 2095.15
$ws.Range("I116").Value = 1167
$ws.Range("K116").Value = 1167
$ws.Range("M116").Value = 1127
$ws.Range("H132").Value = 8599.487999999999
$ws.Range("I132").Value = 5326.4595
$ws.Range("K132").Value = 15979.3785
$ws.Range("M132").Value = -13449.3785
$ws.Range("H136").Value = 2698.5
$ws.Range("I136").Value = 2698.5
$ws.Range("K136").Value = 8095.5
$ws.Range("M136").Value = -5545.5
$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Range("H3").Value = 2095.15
$ws.Range("I3").Value = 1167
$ws.Range("K3").Value = 1167
$ws.Range("M3").Value = -1053
$ws.Range("H86").Value = 2494.6667
$ws.Range("I86").Value = 2075.5
$ws.Range("K86").Value = 2075.5
$ws.Range("M86").Value = -952.5
$ws.Range("H89").Value = 2494.6667
$ws.Range("I89").Value = 2075.5
$ws.Range("K89").Value = 10377.5
$ws.Range("M89").Value = -4761.5
$ws.Range("H94").Value = 3157.842
$ws.Range("I94").Value = 2172.5334
$ws.Range("J94").Value = 6852.75
$ws.Range("K94").Value = 2172.5334
$ws.Range("L94").Value = 6852.75
$ws.Range("M94").Value = -1721.5334
$ws.Range("N94").Value = -7754.75
$ws.Range("H105").Value = 1931.8667
$ws.Range("I105").Value = 1665.0952
$ws.Range("K105").Value = 1665.0952
$ws.Range("M105").Value = 81.90480000000002
$ws.Range("H134").Value = 4572.6924
$ws.Range("I134").Value = 4609.5
$ws.Range("K134").Value = 13828.5
$ws.Range("M134").Value = -11293.5
$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Range("H6").Value = 1874
$ws.Range("I6").Value = 2165.6667
$ws.Range("J6").Value = 999
$ws.Range("K6").Value = 2165.6667
$ws.Range("L6").Value = 999
$ws.Range("M6").Value = -2052.6667
$ws.Range("N6").Value = -1225
$ws.Range("H31").Value = 1769
$ws.Range("I31").Value = 1549.5
$ws.Range("J31").Value = 2082.5715
$ws.Range("K31").Value = 1549.5
$ws.Range("L31").Value = 2082.5715
$ws.Range("M31").Value = -1254.5
$ws.Range("N31").Value = -2672.5715
$ws.Range("H34").Value = 1769
$ws.Range("I34").Value = 1549.5
$ws.Range("J34").Value = 2082.5715
$ws.Range("K34").Value = 1549.5
$ws.Range("L34").Value = 2082.5715
$ws.Range("M34").Value = -1347.5
$ws.Range("N34").Value = -2486.5715
$ws.Range("H50").Value = 30061.334
$ws.Range("I50").Value = 10000
$ws.Range("K50").Value = 10000
$ws.Range("M50").Value = -9375
$ws.Range("H62").Value = 6639.3335
$ws.Range("J62").Value = 6691.3335
$ws.Range("L62").Value = 6691.3335
$ws.Range("N62").Value = -7939.3335
$ws.Range("H65").Value = 6639.3335
$ws.Range("J65").Value = 6691.3335
$ws.Range("L65").Value = 33456.6675
$ws.Range("N65").Value = -39696.6675
$ws.Range("H86").Value = 500001660
$ws.Range("I86").Value = 500001660
$ws.Range("K86").Value = 500001660
$ws.Range("M86").Value = -500000537
$ws.Range("H89").Value = 500001660
$ws.Range("I89").Value = 500001660
$ws.Range("K89").Value = 2500008300
$ws.Range("M89").Value = -2500002684
$ws.Range("H99").Value = 5234.8667
$ws.Range("I99").Value = 4826.4
$ws.Range("J99").Value = 6051.8
$ws.Range("K99").Value = 4826.4
$ws.Range("L99").Value = 6051.8
$ws.Range("M99").Value = -3328.4
$ws.Range("N99").Value = -9047.799999999999
$ws.Range("H122").Value = 1986.7059
$ws.Range("J122").Value = 1987.8
$ws.Range("L122").Value = 5963.4
$ws.Range("N122").Value = -10863.4
$ws.Range("H126").Value = 5234.8667
$ws.Range("I126").Value = 4826.4
$ws.Range("J126").Value = 6051.8
$ws.Range("K126").Value = 14479.2
$ws.Range("L126").Value = 18155.4
$ws.Range("M126").Value = -12009.2
$ws.Range("N126").Value = -23095.4
$ws.Range("H132").Value = 2672.4614
$ws.Range("I132").Value = 2524.3044
$ws.Range("K132").Value = 7572.9132
$ws.Range("M132").Value = -5042.9132
$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Range("H12").Value = 230.34782
$ws.Range("I12").Value = 255
$ws.Range("J12").Value = 192
$ws.Range("K12").Value = 765
$ws.Range("L12").Value = 576
$ws.Range("M12").Value = -592
$ws.Range("N12").Value = -922
$ws.Range("H37").Value = 67158.8
$ws.Range("J37").Value = 67158.8
$ws.Range("L37").Value = 201476.4
$ws.Range("N37").Value = -201700.4
$ws.Range("H39").Value = 3376.6
$ws.Range("J39").Value = 3707.3333
$ws.Range("L39").Value = 11121.9999
$ws.Range("N39").Value = -11709.9999
$ws.Range("H50").Value = 186.3077
$ws.Range("I50").Value = 237.125
$ws.Range("J50").Value = 105
$ws.Range("K50").Value = 711.375
$ws.Range("L50").Value = 315
$ws.Range("M50").Value = -230.375
$ws.Range("N50").Value = -1277
$ws.Range("H53").Value = 186.3077
$ws.Range("I53").Value = 237.125
$ws.Range("J53").Value = 105
$ws.Range("K53").Value = 711.375
$ws.Range("L53").Value = 315
$ws.Range("M53").Value = -230.375
$ws.Range("N53").Value = -1277
$ws.Range("H59").Value = 4999
$ws.Range("I59").Value = 1999
$ws.Range("J59").Value = 7999
$ws.Range("K59").Value = 5997
$ws.Range("L59").Value = 23997
$ws.Range("M59").Value = -5457
$ws.Range("N59").Value = -25077
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = ""
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = ""
$ws.Range("H114").Value = 11617.25
$ws.Range("I114").Value = 555
$ws.Range("J114").Value = 16358.214
$ws.Range("K114").Value = 1665
$ws.Range("L114").Value = 49074.642
$ws.Range("M114").Value = 1589
$ws.Range("N114").Value = -55582.642
$ws.Range("H117").Value = 277470.8
$ws.Range("J117").Value = 277470.8
$ws.Range("L117").Value = 832412.3999999999
$ws.Range("N117").Value = -839296.3999999999
$ws.Range("H121").Value = 1575.2858
$ws.Range("J121").Value = 1843.8823
$ws.Range("L121").Value = 5531.6469
$ws.Range("N121").Value = -8151.6469
$ws.Range("H128").Value = 1499985
$ws.Range("I128").Value = 1499985
$ws.Range("K128").Value = 4499955
$ws.Range("M128").Value = -4494975
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = ""
$ws.Range("N132").Value = ""
$ws = $wb.Worksheets.Item(6)  # GSM
$ws.Range("H102").Value = 4633
$ws.Range("I102").Value = 3499.5
$ws.Range("K102").Value = 3499.5
$ws.Range("M102").Value = -1877.5
$ws.Range("H113").Value = 2142.3
$ws.Range("I113").Value = 2140.6667
$ws.Range("J113").Value = 2144.75
$ws.Range("K113").Value = 2140.6667
$ws.Range("L113").Value = 2144.75
$ws.Range("M113").Value = 29.33329999999978
$ws.Range("N113").Value = -6484.75
$ws.Range("H126").Value = 6166.1
$ws.Range("J126").Value = 8003.5
$ws.Range("L126").Value = 24010.5
$ws.Range("N126").Value = -28950.5
$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Range("H7").Value = 3922.44
$ws.Range("I7").Value = 2936.9333
$ws.Range("J7").Value = 5400.7
$ws.Range("K7").Value = 2936.9333
$ws.Range("L7").Value = 5400.7
$ws.Range("M7").Value = -2824.9333
$ws.Range("N7").Value = -5624.7
$ws.Range("H22").Value = 1070.283
$ws.Range("J22").Value = 1092.6595
$ws.Range("L22").Value = 1092.6595
$ws.Range("N22").Value = -1682.6595
$ws.Range("H27").Value = 1070.283
$ws.Range("J27").Value = 1092.6595
$ws.Range("L27").Value = 1092.6595
$ws.Range("N27").Value = -1306.6595
$ws.Range("H40").Value = 7415.2
$ws.Range("I40").Value = 6652
$ws.Range("K40").Value = 6652
$ws.Range("M40").Value = -6516
$ws.Range("H41").Value = 39635.285
$ws.Range("I41").Value = 34974
$ws.Range("J41").Value = 41499.8
$ws.Range("K41").Value = 34974
$ws.Range("L41").Value = 41499.8
$ws.Range("M41").Value = -34536
$ws.Range("N41").Value = -42375.8
$ws.Range("H42").Value = 21198
$ws.Range("I42").Value = 21634
$ws.Range("J42").Value = 19890
$ws.Range("K42").Value = 21634
$ws.Range("L42").Value = 19890
$ws.Range("M42").Value = -21071
$ws.Range("N42").Value = -21016
$ws.Range("H46").Value = 1574.75
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1574.75
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 1574.75
$ws.Range("M46").Value = ""
$ws.Range("N46").Value = -1950.75
$ws.Range("H49").Value = 21198
$ws.Range("I49").Value = 21634
$ws.Range("J49").Value = 19890
$ws.Range("K49").Value = 21634
$ws.Range("L49").Value = 19890
$ws.Range("M49").Value = -21487
$ws.Range("N49").Value = -20184
$ws.Range("H55").Value = 1137.4546
$ws.Range("I55").Value = 559.4286
$ws.Range("J55").Value = 2149
$ws.Range("K55").Value = 559.4286
$ws.Range("L55").Value = 2149
$ws.Range("M55").Value = -386.4286
$ws.Range("N55").Value = -2495
$ws.Range("H68").Value = 3532.4
$ws.Range("I68").Value = 1340.6364
$ws.Range("J68").Value = 9559.75
$ws.Range("K68").Value = 1340.6364
$ws.Range("L68").Value = 9559.75
$ws.Range("M68").Value = -591.6364000000001
$ws.Range("N68").Value = -11057.75
$ws.Range("H71").Value = 3532.4
$ws.Range("I71").Value = 1340.6364
$ws.Range("J71").Value = 9559.75
$ws.Range("K71").Value = 6703.182000000001
$ws.Range("L71").Value = 47798.75
$ws.Range("M71").Value = -2959.182000000001
$ws.Range("N71").Value = -55286.75
$ws.Range("H82").Value = 2435.8572
$ws.Range("J82").Value = 4623.6665
$ws.Range("L82").Value = 4623.6665
$ws.Range("N82").Value = -5345.6665
$ws.Range("H85").Value = 2435.8572
$ws.Range("J85").Value = 4623.6665
$ws.Range("L85").Value = 4623.6665
$ws.Range("N85").Value = -7119.6665
$ws.Range("H93").Value = 5007.4
$ws.Range("I93").Value = 2969.875
$ws.Range("K93").Value = 2969.875
$ws.Range("M93").Value = -1721.875
$ws.Range("H94").Value = 80000
$ws.Range("J94").Value = 80000
$ws.Range("L94").Value = 80000
$ws.Range("N94").Value = -81352
$ws.Range("H121").Value = 200000
$ws.Range("J121").Value = 200000
$ws.Range("L121").Value = 200000
$ws.Range("N121").Value = -203494
$ws.Range("H122").Value = 15202
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 15202
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 45606
$ws.Range("M122").Value = ""
$ws.Range("N122").Value = -50506
$ws.Range("H126").Value = 3922.44
$ws.Range("I126").Value = 2936.9333
$ws.Range("J126").Value = 5400.7
$ws.Range("K126").Value = 8810.7999
$ws.Range("L126").Value = 16202.1
$ws.Range("M126").Value = -6340.7999
$ws.Range("N126").Value = -21142.1
$ws.Range("H132").Value = 5202.778
$ws.Range("I132").Value = 5202.778
$ws.Range("K132").Value = 15608.334
$ws.Range("M132").Value = -13078.334
$ws.Range("H136").Value = 6686.273
$ws.Range("I136").Value = 5959.8887
$ws.Range("J136").Value = 9955
$ws.Range("K136").Value = 17879.6661
$ws.Range("L136").Value = 29865
$ws.Range("M136").Value = -15329.6661
$ws.Range("N136").Value = -34965
$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Range("H12").Value = 6996.2
$ws.Range("I12").Value = 4995.25
$ws.Range("J12").Value = 15000
$ws.Range("K12").Value = 4995.25
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = -4853.25
$ws.Range("N12").Value = -15284
$ws.Range("H37").Value = 48762.5
$ws.Range("I37").Value = 48762.5
$ws.Range("K37").Value = 48762.5
$ws.Range("M37").Value = -48559.5
$ws.Range("H62").Value = 55575.434
$ws.Range("I62").Value = 81739.8
$ws.Range("K62").Value = 81739.8
$ws.Range("M62").Value = -81115.8
$ws.Range("H65").Value = 55575.434
$ws.Range("I65").Value = 81739.8
$ws.Range("K65").Value = 408699
$ws.Range("M65").Value = -405579
$ws.Range("H107").Value = 586.6667
$ws.Range("I107").Value = 358.33334
$ws.Range("K107").Value = 1075.00002
$ws.Range("M107").Value = 844.9999800000001
$ws.Range("H122").Value = 4774.7417
$ws.Range("J122").Value = 3286.6924
$ws.Range("L122").Value = 9860.0772
$ws.Range("N122").Value = -14760.0772
$ws.Range("H126").Value = 3109
$ws.Range("I126").Value = 2603.3333
$ws.Range("K126").Value = 7809.999899999999
$ws.Range("M126").Value = -5339.999899999999
$ws.Range("H136").Value = 8989.4
$ws.Range("I136").Value = 10699.286
$ws.Range("K136").Value = 32097.858
$ws.Range("M136").Value = -29547.858
